$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1288.1666
$ws.Range("I15").Value = 1288.1666
$ws.Range("K15").Value = 3864.4998
$ws.Range("M15").Value = -3695.4998
$ws.Range("H16").Value = 5000
$ws.Range("I16").Value = 5000
$ws.Range("K16").Value = 5000
$ws.Range("M16").Value = -4770
$ws.Range("H62").Value = 7175.8
$ws.Range("I62").Value = 1974.75
$ws.Range("J62").Value = 27980
$ws.Range("K62").Value = 1974.75
$ws.Range("L62").Value = 27980
$ws.Range("M62").Value = -1350.75
$ws.Range("N62").Value = -29228
$ws.Range("H65").Value = 7175.8
$ws.Range("I65").Value = 1974.75
$ws.Range("J65").Value = 27980
$ws.Range("K65").Value = 9873.75
$ws.Range("L65").Value = 139900
$ws.Range("M65").Value = -6753.75
$ws.Range("N65").Value = -146140
$ws.Range("H98").Value = 2653.3333
$ws.Range("I98").Value = 2783.3333
$ws.Range("K98").Value = 2783.3333
$ws.Range("M98").Value = -1285.3333
$ws.Range("H106").Value = 1033.875
$ws.Range("I106").Value = 895.8571
$ws.Range("K106").Value = 895.8571
$ws.Range("M106").Value = -264.8570999999999
$ws.Range("H116").Value = 15875.625
$ws.Range("I116").Value = 27751.25
$ws.Range("J116").Value = 4000
$ws.Range("K116").Value = 27751.25
$ws.Range("L116").Value = 4000
$ws.Range("M116").Value = -24309.25
$ws.Range("N116").Value = -10884
$ws.Range("H122").Value = 2653.3333
$ws.Range("I122").Value = 2783.3333
$ws.Range("K122").Value = 8349.999899999999
$ws.Range("M122").Value = -5899.999899999999
$ws.Range("H132").Value = 1023.2826
$ws.Range("I132").Value = 1070.6757
$ws.Range("J132").Value = 828.44446
$ws.Range("K132").Value = 3212.0271
$ws.Range("L132").Value = 2485.33338
$ws.Range("M132").Value = -682.0271000000002
$ws.Range("N132").Value = -7545.33338
$ws.Range("H135").Value = 499.2353
$ws.Range("I135").Value = 440.57144
$ws.Range("K135").Value = 3965.14296
$ws.Range("M135").Value = -1430.14296
$ws.Range("H137").Value = 54067.58
$ws.Range("I137").Value = 776.8889
$ws.Range("J137").Value = 102029.2
$ws.Range("K137").Value = 2330.6667
$ws.Range("L137").Value = 306087.6
$ws.Range("M137").Value = 219.3332999999998
$ws.Range("N137").Value = -311187.6
$ws.Range("H138").Value = 3289.2263
$ws.Range("J138").Value = 2856.658
$ws.Range("L138").Value = 8569.974
$ws.Range("N138").Value = -18849.974
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2761.6748
$ws.Range("I32").Value = 2197.6487
$ws.Range("J32").Value = 7399.222
$ws.Range("K32").Value = 2197.6487
$ws.Range("L32").Value = 7399.222
$ws.Range("M32").Value = -1910.6487
$ws.Range("N32").Value = -7973.222
$ws.Range("H37").Value = 14000
$ws.Range("J37").Value = 14000
$ws.Range("L37").Value = 14000
$ws.Range("N37").Value = -14546
$ws.Range("H44").Value = 8500
$ws.Range("I44").Value = 8500
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 8500
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("M44").Value = -8012
$ws.Range("H45").Value = 2161.4
$ws.Range("I45").Value = 1003.6667
$ws.Range("J45").Value = 3898
$ws.Range("K45").Value = 1003.6667
$ws.Range("L45").Value = 3898
$ws.Range("M45").Value = -626.6667
$ws.Range("N45").Value = -4652
$ws.Range("H61").Value = 5157.857
$ws.Range("I61").Value = 1995
$ws.Range("J61").Value = 5685
$ws.Range("K61").Value = 1995
$ws.Range("L61").Value = 5685
$ws.Range("M61").Value = -1783
$ws.Range("N61").Value = -6109
$ws.Range("H74").Value = 839.36
$ws.Range("I74").Value = 773.2174
$ws.Range("J74").Value = 1600
$ws.Range("K74").Value = 773.2174
$ws.Range("L74").Value = 1600
$ws.Range("M74").Value = 100.7826
$ws.Range("N74").Value = -3348
$ws.Range("H77").Value = 839.36
$ws.Range("I77").Value = 773.2174
$ws.Range("J77").Value = 1600
$ws.Range("K77").Value = 3866.087
$ws.Range("L77").Value = 8000
$ws.Range("M77").Value = 501.913
$ws.Range("N77").Value = -16736
$ws.Range("H122").Value = 1534.6923
$ws.Range("I122").Value = 1566.75
$ws.Range("J122").Value = 1150
$ws.Range("K122").Value = 4700.25
$ws.Range("L122").Value = 3450
$ws.Range("M122").Value = -2250.25
$ws.Range("N122").Value = -8350
$ws.Range("H136").Value = 5157.857
$ws.Range("I136").Value = 1995
$ws.Range("J136").Value = 5685
$ws.Range("K136").Value = 5985
$ws.Range("L136").Value = 17055
$ws.Range("M136").Value = -3435
$ws.Range("N136").Value = -22155
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2095.6
$ws.Range("I20").Value = 1793.6666
$ws.Range("K20").Value = 1793.6666
$ws.Range("M20").Value = -1546.6666
$ws.Range("H107").Value = 1840.1333
$ws.Range("J107").Value = 3330.6667
$ws.Range("L107").Value = 3330.6667
$ws.Range("N107").Value = -7170.6667
$ws.Range("H134").Value = 6613.4707
$ws.Range("I134").Value = 6828.6
$ws.Range("K134").Value = 20485.8
$ws.Range("M134").Value = -17950.8
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2171.8809
$ws.Range("I31").Value = 1200
$ws.Range("J31").Value = 2769.9614
$ws.Range("K31").Value = 1200
$ws.Range("L31").Value = 2769.9614
$ws.Range("M31").Value = -905
$ws.Range("N31").Value = -3359.9614
$ws.Range("H34").Value = 2171.8809
$ws.Range("I34").Value = 1200
$ws.Range("J34").Value = 2769.9614
$ws.Range("K34").Value = 1200
$ws.Range("L34").Value = 2769.9614
$ws.Range("M34").Value = -998
$ws.Range("N34").Value = -3173.9614
$ws.Range("H105").Value = 926.7143
$ws.Range("I105").Value = 926.7143
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 926.7143
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 820.2857
$ws.Range("N105").ClearContents()
$ws.Range("H132").Value = 2973.52
$ws.Range("I132").Value = 949
$ws.Range("J132").Value = 3926.2354
$ws.Range("K132").Value = 2847
$ws.Range("L132").Value = 11778.7062
$ws.Range("M132").Value = -317
$ws.Range("N132").Value = -16838.7062
$ws.Range("H134").Value = 2306.16
$ws.Range("I134").Value = 1171.2632
$ws.Range("K134").Value = 3513.7896
$ws.Range("M134").Value = -978.7896000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1607244
$ws.Range("I4").Value = 1577032
$ws.Range("K4").Value = 4731096
$ws.Range("M4").Value = -4730984
$ws.Range("H16").Value = 333333340
$ws.Range("I16").Value = 333333340
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1000000020
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -999999847
$ws.Range("N16").ClearContents()
$ws.Range("H56").Value = 916720.94
$ws.Range("I56").Value = 916720.94
$ws.Range("K56").Value = 916720.94
$ws.Range("M56").Value = -916190.94
$ws.Range("H68").Value = 2204.8958
$ws.Range("J68").Value = 2983.7144
$ws.Range("L68").Value = 8951.143199999999
$ws.Range("N68").Value = -10573.1432
$ws.Range("H71").Value = 2204.8958
$ws.Range("J71").Value = 2983.7144
$ws.Range("L71").Value = 26853.4296
$ws.Range("N71").Value = -34965.4296
$ws.Range("H107").Value = 2354.4595
$ws.Range("J107").Value = 2425.3
$ws.Range("L107").Value = 7275.900000000001
$ws.Range("N107").Value = -11115.9
$ws.Range("H113").Value = 68012.734
$ws.Range("I113").Value = 202357.2
$ws.Range("J113").Value = 840.5
$ws.Range("K113").Value = 607071.6000000001
$ws.Range("L113").Value = 2521.5
$ws.Range("M113").Value = -604901.6000000001
$ws.Range("N113").Value = -6861.5
$ws.Range("H131").Value = 10219489
$ws.Range("J131").Value = 17100.092
$ws.Range("L131").Value = 51300.276
$ws.Range("N131").Value = -61380.276
$ws.Range("H133").Value = 50001930
$ws.Range("I133").Value = 50001930
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 150005790
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -150000730
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 3999.25
$ws.Range("I136").Value = 3750
$ws.Range("J136").Value = 4248.5
$ws.Range("K136").Value = 11250
$ws.Range("L136").Value = 12745.5
$ws.Range("M136").Value = -6150
$ws.Range("N136").Value = -22945.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 6056.5557
$ws.Range("I32").Value = 5579.875
$ws.Range("J32").Value = 9870
$ws.Range("K32").Value = 5579.875
$ws.Range("L32").Value = 9870
$ws.Range("M32").Value = -5262.875
$ws.Range("N32").Value = -10504
$ws.Range("H122").Value = 5338
$ws.Range("I122").Value = 3903.3635
$ws.Range("J122").Value = 6772.636
$ws.Range("K122").Value = 11710.0905
$ws.Range("L122").Value = 20317.908
$ws.Range("M122").Value = -9260.0905
$ws.Range("N122").Value = -25217.908
$ws.Range("H136").Value = 4755.647
$ws.Range("I136").Value = 2682.2856
$ws.Range("J136").Value = 6207
$ws.Range("K136").Value = 8046.8568
$ws.Range("L136").Value = 18621
$ws.Range("M136").Value = -5496.8568
$ws.Range("N136").Value = -23721
